# Add 2022-Q4 data
# ------------------------------------------------------------------
# 1) Insert a brand-new worksheet "2022-Q4" right before "2022-Q3"
#    (i.e. right after the "总计" summary sheet) and fill it with the
#    per-fund holdings table for the new quarter.
# 2) Update the "总计" (summary) sheet: a new row is inserted at the
#    top of the data (row 2) holding the 2022-Q4 totals, and every
#    later row's B/C/D (date/count/value) slide down by one slot, with
#    a brand-new trailing row appended for what used to be the last
#    row. Column A is just the static 0-based rank index, so it is
#    left alone except for the newly-appended row.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Helper: write a value into a cell as *text* (shared-string / inlineStr
# semantics) even when it looks like a number, mirroring how the source
# data was authored, then drop the left-over "@" number-format style so
# the cell ends up on the default (un-styled) cell format - exactly like
# its sibling cells in the existing sheets.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Helper: copy the *formatting only* of $srcRange onto $dstRange (style,
# borders, font, alignment - no value).
function Copy-Format($srcRange, $dstRange) {
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
}

# ===================================================================
# Step 1 - brand-new "2022-Q4" worksheet
# ===================================================================

$summary = $wb.Worksheets.Item("总计")
$anchor = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# Header row (B1:H1) - bold/bordered style (same as every other quarter
# sheet), copied from the summary sheet's header cell. These header
# labels are plain (non-numeric-looking) text, so a normal .Value
# assignment already stores them as text - no need for the
# NumberFormat/ClearFormats dance (which would strip the pasted style
# back off again).
Copy-Format $summary.Range("B1") $q4.Range("B1:H1")

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Data rows 2..10 - column A is the bold/bordered rank index, B-G are
# plain text, H is a plain number.
$rows = @(
    @(0, "004932", "招商丰拓灵活配置混合A", "5.03",  "90.46", "5.70", "0.2867", 3),
    @(1, "004933", "招商丰拓灵活配置混合C", "4.90",  "90.46", "5.70", "0.2793", 3),
    @(2, "002657", "招商安裕灵活配置混合A", "13.66", "33.20", "2.03", "0.2773", 6),
    @(3, "002658", "招商安裕灵活配置混合A", "4.26",  "33.20", "2.03", "0.0865", 6),
    @(4, "004143", "招商盛合灵活配置混合C", "1.92",  "55.65", "4.28", "0.0822", 5),
    @(5, "002581", "招商丰凯灵活配置混合A", "1.68",  "37.13", "2.10", "0.0353", 3),
    @(6, "015206", "招商安裕灵活配置混合D", "1.25",  "33.20", "2.03", "0.0254", 6),
    @(7, "002582", "招商丰凯灵活配置混合C", "0.96",  "37.13", "2.10", "0.0202", 3),
    @(8, "004142", "招商盛合灵活配置混合A", "0.07",  "55.65", "4.28", "0.0030", 5)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    Copy-Format $summary.Range("A2") $q4.Range("A$r")
    $q4.Range("A$r").Value = $row[0]

    Set-TextValue $q4.Range("B$r") $row[1]
    Set-TextValue $q4.Range("C$r") $row[2]
    Set-TextValue $q4.Range("D$r") $row[3]
    Set-TextValue $q4.Range("E$r") $row[4]
    Set-TextValue $q4.Range("F$r") $row[5]
    Set-TextValue $q4.Range("G$r") $row[6]

    $q4.Range("H$r").Value = $row[7]
}

# ===================================================================
# Step 2 - update the "总计" sheet
# ===================================================================

# Remember the current (pre-edit) B/C/D content for rows 2..8 so we can
# slide it down one row before overwriting row 2 with the new quarter.
$oldB = @{}
$oldC = @{}
$oldD = @{}
for ($r = 2; $r -le 8; $r++) {
    $oldB[$r] = $summary.Range("B$r").Value2
    $oldC[$r] = $summary.Range("C$r").Value2
    $oldD[$r] = $summary.Range("D$r").Value2
}

# Slide rows 8->9, 7->8, ..., 2->3 (walk from the bottom up so we never
# overwrite a value before it has been read/copied).
for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    Set-TextValue $summary.Range("B$dest") $oldB[$r]
    $summary.Range("C$dest").Value = $oldC[$r]
    $summary.Range("D$dest").Value = $oldD[$r]
}

# New row 9's rank index + style (copy from row 8, which already carries
# the bold/bordered "s=2" look).
Copy-Format $summary.Range("A8") $summary.Range("A9")
$summary.Range("A9").Value = 7

# Finally, drop the brand-new 2022-Q4 totals into row 2.
Set-TextValue $summary.Range("B2") "2022-Q4"
$summary.Range("C2").Value = 9
$summary.Range("D2").Value = 1.1

# Restore the original active tab (inserting a sheet makes it active,
# but the last sheet - "2020-Q4" - was the one actually selected before
# our edit) so the view state is left exactly as it was found.
$wb.Worksheets.Item("2020-Q4").Activate()
